$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# ---- ALC ----
# row 4
$ws1.Range("H4").Value = 465.33334
$ws1.Range("I4").Value = 0
$ws1.Range("K4").Value = 0
$ws1.Range("M4").ClearContents()
# row 17
$ws1.Range("H17").Value = 941.8039
$ws1.Range("J17").Value = 940.375
$ws1.Range("L17").Value = 2821.125
$ws1.Range("N17").Value = -3157.125
# row 28
$ws1.Range("H28").Value = 605.5
$ws1.Range("I28").Value = 395.8421
$ws1.Range("K28").Value = 395.8421
$ws1.Range("M28").Value = 89.15789999999998
# row 86
$ws1.Range("H86").Value = 3512.75
$ws1.Range("I86").Value = 1052
$ws1.Range("K86").Value = 1052
$ws1.Range("M86").Value = 71
# row 89
$ws1.Range("H89").Value = 3512.75
$ws1.Range("I89").Value = 1052
$ws1.Range("K89").Value = 5260
$ws1.Range("M89").Value = 356
# row 111
$ws1.Range("H111").Value = 4104
$ws1.Range("I111").Value = 3125.5
$ws1.Range("K111").Value = 9376.5
$ws1.Range("M111").Value = -6309.5
# row 116
$ws1.Range("H116").Value = 206527.58
$ws1.Range("I116").Value = 65283.2
$ws1.Range("J116").Value = 334931.53
$ws1.Range("K116").Value = 65283.2
$ws1.Range("L116").Value = 334931.53
$ws1.Range("M116").Value = -61841.2
$ws1.Range("N116").Value = -341815.53
# row 141
$ws1.Range("H141").Value = 1555.75
$ws1.Range("I141").Value = 1555.75
$ws1.Range("K141").Value = 4667.25
$ws1.Range("M141").Value = 512.75

# ---- ARM ----
# row 32
$ws2.Range("H32").Value = 6329931.5
$ws2.Range("I32").Value = 6411073
$ws2.Range("K32").Value = 6411073
$ws2.Range("M32").Value = -6410786
# row 45
$ws2.Range("H45").Value = 2018.1428
$ws2.Range("I45").Value = 1983.8
$ws2.Range("K45").Value = 1983.8
$ws2.Range("M45").Value = -1606.8
# row 97
$ws2.Range("H97").Value = 1116.409
$ws2.Range("I97").Value = 1116.25
$ws2.Range("K97").Value = 1116.25
$ws2.Range("M97").Value = -620.25
# row 102
$ws2.Range("H102").Value = 61624
$ws2.Range("I102").Value = 72749
$ws2.Range("J102").Value = 5999
$ws2.Range("K102").Value = 72749
$ws2.Range("L102").Value = 5999
$ws2.Range("M102").Value = -71127
$ws2.Range("N102").Value = -9243
# row 110
$ws2.Range("H110").Value = 2436
$ws2.Range("I110").Value = 3800
$ws2.Range("J110").Value = 1072
$ws2.Range("K110").Value = 3800
$ws2.Range("L110").Value = 1072
$ws2.Range("M110").Value = -1755
$ws2.Range("N110").Value = -5162
# row 132
$ws2.Range("H132").Value = 379204.72
$ws2.Range("I132").Value = 476451.25
$ws2.Range("K132").Value = 1429353.75
$ws2.Range("M132").Value = -1426823.75

# ---- BSM ----
# row 107
$ws3.Range("H107").Value = 5114.9287
$ws3.Range("I107").Value = 5175.75
$ws3.Range("K107").Value = 5175.75
$ws3.Range("M107").Value = -3255.75
# row 134
$ws3.Range("H134").Value = 1309969.9
$ws3.Range("I134").Value = 1513430.5
$ws3.Range("K134").Value = 4540291.5
$ws3.Range("M134").Value = -4537756.5

# ---- CRP ----
# row 31
$ws4.Range("H31").Value = 142155.31
$ws4.Range("I31").Value = 392064
$ws4.Range("J31").Value = 24551.234
$ws4.Range("K31").Value = 392064
$ws4.Range("L31").Value = 24551.234
$ws4.Range("M31").Value = -391769
$ws4.Range("N31").Value = -25141.234
# row 34
$ws4.Range("H34").Value = 142155.31
$ws4.Range("I34").Value = 392064
$ws4.Range("J34").Value = 24551.234
$ws4.Range("K34").Value = 392064
$ws4.Range("L34").Value = 24551.234
$ws4.Range("M34").Value = -391862
$ws4.Range("N34").Value = -24955.234
# row 39
$ws4.Range("H39").Value = 1000
$ws4.Range("I39").Value = 1000
$ws4.Range("K39").Value = 1000
$ws4.Range("M39").Value = -609
# row 49
$ws4.Range("H49").Value = 1000
$ws4.Range("I49").Value = 1000
$ws4.Range("K49").Value = 1000
$ws4.Range("M49").Value = -818
# row 53
$ws4.Range("H53").Value = 55950
$ws4.Range("J53").Value = 55950
$ws4.Range("L53").Value = 55950
$ws4.Range("N53").Value = -57164
# row 58
$ws4.Range("H58").Value = 366265.3
$ws4.Range("I58").Value = 538047.2
$ws4.Range("J58").Value = 7085.091
$ws4.Range("K58").Value = 538047.2
$ws4.Range("L58").Value = 7085.091
$ws4.Range("M58").Value = -537844.2
$ws4.Range("N58").Value = -7491.091
# row 88
$ws4.Range("H88").Value = 46989.668
$ws4.Range("J88").Value = 47787.6
$ws4.Range("L88").Value = 47787.6
$ws4.Range("N88").Value = -48599.6
# row 91
$ws4.Range("H91").Value = 46989.668
$ws4.Range("J91").Value = 47787.6
$ws4.Range("L91").Value = 47787.6
$ws4.Range("N91").Value = -50595.6
# row 99
$ws4.Range("H99").Value = 2223.3076
$ws4.Range("I99").Value = 1986.8572
$ws4.Range("J99").Value = 2499.1667
$ws4.Range("K99").Value = 1986.8572
$ws4.Range("L99").Value = 2499.1667
$ws4.Range("M99").Value = -488.8571999999999
$ws4.Range("N99").Value = -5495.1667
# row 103
$ws4.Range("H103").Value = 14690.833
$ws4.Range("I103").Value = 14690.833
$ws4.Range("K103").Value = 14690.833
$ws4.Range("M103").Value = -13518.833
# row 107
$ws4.Range("H107").Value = 1332.3334
$ws4.Range("I107").Value = 1300
$ws4.Range("J107").Value = 1348.5
$ws4.Range("K107").Value = 1300
$ws4.Range("L107").Value = 1348.5
$ws4.Range("M107").Value = 620
$ws4.Range("N107").Value = -5188.5
# row 111
$ws4.Range("H111").Value = 100702
$ws4.Range("J111").Value = 100702
$ws4.Range("L111").Value = 100702
$ws4.Range("N111").Value = -108882
# row 118
$ws4.Range("H118").Value = 39871
$ws4.Range("J118").Value = 39871
$ws4.Range("L118").Value = 39871
$ws4.Range("N118").Value = -43185
# row 126
$ws4.Range("H126").Value = 2223.3076
$ws4.Range("I126").Value = 1986.8572
$ws4.Range("J126").Value = 2499.1667
$ws4.Range("K126").Value = 5960.571599999999
$ws4.Range("L126").Value = 7497.500100000001
$ws4.Range("M126").Value = -3490.571599999999
$ws4.Range("N126").Value = -12437.5001
# row 134
$ws4.Range("H134").Value = 16019.542
$ws4.Range("I134").Value = 19203.316
$ws4.Range("J134").Value = 3921.2
$ws4.Range("K134").Value = 57609.948
$ws4.Range("L134").Value = 11763.6
$ws4.Range("M134").Value = -55074.948
$ws4.Range("N134").Value = -16833.6
# row 136
$ws4.Range("H136").Value = 366265.3
$ws4.Range("I136").Value = 538047.2
$ws4.Range("J136").Value = 7085.091
$ws4.Range("K136").Value = 1614141.6
$ws4.Range("L136").Value = 21255.273
$ws4.Range("M136").Value = -1611591.6
$ws4.Range("N136").Value = -26355.273

# ---- CUL ----
# row 4
$ws5.Range("H4").Value = 36604680
$ws5.Range("I4").Value = 38980560
$ws5.Range("J4").Value = 25200440
$ws5.Range("K4").Value = 116941680
$ws5.Range("L4").Value = 75601320
$ws5.Range("M4").Value = -116941568
$ws5.Range("N4").Value = -75601544
# row 7
$ws5.Range("H7").Value = 375099.75
$ws5.Range("I7").Value = 199
$ws5.Range("K7").Value = 597
$ws5.Range("M7").Value = -485
# row 133
$ws5.Range("H133").Value = 7691.1665
$ws5.Range("I133").Value = 5636
$ws5.Range("K133").Value = 16908
$ws5.Range("M133").Value = -11848

# ---- GSM ----
# row 7
$ws6.Range("H7").Value = 30000
$ws6.Range("I7").Value = 30000
$ws6.Range("K7").Value = 30000
$ws6.Range("M7").Value = -29888
# row 8
$ws6.Range("H8").Value = 30000
$ws6.Range("I8").Value = 30000
$ws6.Range("K8").Value = 30000
$ws6.Range("M8").Value = -29861
# row 14
$ws6.Range("H14").Value = 1855145.9
$ws6.Range("I14").Value = 2040160.4
$ws6.Range("J14").Value = 5000
$ws6.Range("K14").Value = 2040160.4
$ws6.Range("L14").Value = 5000
$ws6.Range("M14").Value = -2039992.4
$ws6.Range("N14").Value = -5336
# row 36
$ws6.Range("H36").Value = 15000
$ws6.Range("I36").Value = 15000
$ws6.Range("J36").Value = 0
$ws6.Range("K36").Value = 15000
$ws6.Range("L36").Value = 0
$ws6.Range("M36").Value = -14515
$ws6.Range("N36").ClearContents()
# row 43
$ws6.Range("H43").Value = 5084.5
$ws6.Range("I43").Value = 1356
$ws6.Range("J43").Value = 19998.5
$ws6.Range("K43").Value = 1356
$ws6.Range("L43").Value = 19998.5
$ws6.Range("M43").Value = -1205
$ws6.Range("N43").Value = -20300.5
# row 48
$ws6.Range("H48").Value = 0
$ws6.Range("I48").Value = 0
$ws6.Range("K48").Value = 0
$ws6.Range("M48").ClearContents()
# row 107
$ws6.Range("H107").Value = 52457.43
$ws6.Range("I107").Value = 118499
$ws6.Range("K107").Value = 118499
$ws6.Range("M107").Value = -116579
# row 113
$ws6.Range("H113").Value = 4665
$ws6.Range("I113").Value = 4665
$ws6.Range("J113").Value = 0
$ws6.Range("K113").Value = 4665
$ws6.Range("L113").Value = 0
$ws6.Range("M113").Value = -2495
$ws6.Range("N113").ClearContents()
# row 129
$ws6.Range("H129").Value = 0
$ws6.Range("J129").Value = 0
$ws6.Range("L129").Value = 0
$ws6.Range("N129").ClearContents()
# row 132
$ws6.Range("H132").Value = 19107034
$ws6.Range("I132").Value = 28928174
$ws6.Range("J132").Value = 10372.056
$ws6.Range("K132").Value = 86784522
$ws6.Range("L132").Value = 31116.168
$ws6.Range("M132").Value = -86781992
$ws6.Range("N132").Value = -36176.16800000001

# ---- LTW ----
# row 68
$ws7.Range("H68").Value = 3299.75
$ws7.Range("J68").Value = 3499.5
$ws7.Range("L68").Value = 3499.5
$ws7.Range("N68").Value = -4997.5
# row 71
$ws7.Range("H71").Value = 3299.75
$ws7.Range("J71").Value = 3499.5
$ws7.Range("L71").Value = 17497.5
$ws7.Range("N71").Value = -24985.5

# ---- WVR ----
# row 82
$ws8.Range("H82").Value = 36993.75
$ws8.Range("I82").Value = 0
$ws8.Range("J82").Value = 36993.75
$ws8.Range("K82").Value = 0
$ws8.Range("L82").Value = 36993.75
$ws8.Range("M82").ClearContents()
$ws8.Range("N82").Value = -37759.75
# row 85
$ws8.Range("H85").Value = 36993.75
$ws8.Range("I85").Value = 0
$ws8.Range("J85").Value = 36993.75
$ws8.Range("K85").Value = 0
$ws8.Range("L85").Value = 36993.75
$ws8.Range("M85").ClearContents()
$ws8.Range("N85").Value = -39645.75
# row 132
$ws8.Range("H132").Value = 7298006
$ws8.Range("I132").Value = 9555079
$ws8.Range("K132").Value = 28665237
$ws8.Range("M132").Value = -28662707
# row 136
$ws8.Range("H136").Value = 12702114
$ws8.Range("I136").Value = 14153710
$ws8.Range("K136").Value = 42461130
$ws8.Range("M136").Value = -42458580
